$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.345780849456787
$ws.Range("B1").Value = 3.303092241287231
$ws.Range("C1").Value = 5.47325325012207
$ws.Range("D1").Value = 1.687505602836609
$ws.Range("E1").Value = 0.9854691624641418
